# Update the embedded build timestamp throughout the workbook.
# Old build stamp -> new build stamp (per the "mines - January 30" release diff).

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

# --- "About" sheet: A2 (version banner) and A6 (recommended citation) ---
$aboutWs = $wb.Worksheets.Item("About")

$a2 = $aboutWs.Range("A2")
$a2.Value2 = $a2.Value2.Replace($oldStamp, $newStamp)

$a6 = $aboutWs.Range("A6")
$a6.Value2 = $a6.Value2.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet: S2:S13 (build_version column) ---
$boundWs = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 13; $row++) {
    $cell = $boundWs.Cells.Item($row, 19)   # column S = 19
    $cell.Value2 = $cell.Value2.Replace($oldStamp, $newStamp)
}
